# Edits are issued in the same order the original author made them so that
# newly-introduced shared strings land at the same table positions Excel
# would naturally produce (new unique strings are appended as they are
# first used).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Seguimiento")

# --- Row 9 ---
# F9 status text changes ("En revisión de editor" -> "Revisado - Pnte cambio motores M")
$ws.Range("F9").Value = "Revisado - Pnte cambio motores M"
# Row height grows to fit the now-longer wrapped text
$ws.Rows.Item(9).RowHeight = 30.75

# --- Row 10 ---
# B10 / C10 get date values (serial 42107 = 2015-04-13)
$ws.Range("B10").Value = [DateTime]::FromOADate(42107)
$ws.Range("C10").Value = [DateTime]::FromOADate(42107)
# C10 picks up B10's formatting (medium-left-border date style), matching
# the look already used for other fully-dated rows (e.g. row 12, 18)
$ws.Range("B10").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$excel.CutCopyMode = 0
# F10 status text changes ("En revisión de editor" -> "En Revisión")
$ws.Range("F10").Value = "En Revisión"

# --- Row 6 ---
# B6 gets a date value (serial 42110 = 2015-04-16)
$ws.Range("B6").Value = [DateTime]::FromOADate(42110)
# F6 status text changes ("Cuaderno de estudio" -> "En revisión")
$ws.Range("F6").Value = "En revisión"

# --- Row 7 ---
# F7 status text is cleared (was "Cuaderno de estudio")
$ws.Range("F7").ClearContents()

# --- View / selection ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("B6").Select()
